$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column H (8) from 13.0 to 14.5 (stored OOXML "width" units).
# ColumnWidth is expressed in character-width units and gets re-derived
# from pixels on save, so 13.67 is the value that round-trips to 14.5.
$ws.Columns.Item(8).ColumnWidth = 13.67

# Move the "nome_objeto" value from column H into column G for rows 2-11,
# and only keep/set "sim" in column H where it still applies; otherwise
# fully clear the H cell so it no longer appears in the sheet.
$ws.Range("G2").Value2 = "gavetas_fechadas"
$ws.Range("H2").Value2 = "sim"

$ws.Range("G3").Value2 = "gavetas_abertas"
$ws.Range("H3").Clear()

$ws.Range("G4").Value2 = "frente_aberta"
$ws.Range("H4").Clear()

$ws.Range("G5").Value2 = "frente_fechada"
$ws.Range("H5").Value2 = "sim"

$ws.Range("G6").Value2 = "estrutura"
$ws.Range("H6").Value2 = "sim"

$ws.Range("G7").Value2 = "gavetas_fechadas"
$ws.Range("H7").Clear()

$ws.Range("G8").Value2 = "gavetas_abertas"
$ws.Range("H8").Clear()

$ws.Range("G9").Value2 = "frente_aberta"
$ws.Range("H9").Clear()

$ws.Range("G10").Value2 = "frente_fechada"
$ws.Range("H10").Clear()

$ws.Range("G11").Value2 = "estrutura"
$ws.Range("H11").Clear()

# These rows keep their existing G value, but drop the now-redundant "sim" in H.
$ws.Range("H20").Clear()
$ws.Range("H28").Clear()
$ws.Range("H29").Clear()
$ws.Range("H30").Clear()
